$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet3"

$ws33 = $wb.Worksheets.Item("Sheet33")
$ws33.Name = "1"

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("5").Insert()
$ws1.Range("A4:H4").Copy()
$ws1.Range("A5").PasteSpecial()

[void]$ws1.Activate()
[void]$ws1.Range("E14").Select()
